$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / title to reflect the new date
$ws.Name = "Through 2021-09-20"

# Update the September row label (through date)
$ws.Range("A10").Value = "September (through 09-20)"

# Update September (row 10) counts per year (columns B..H = 2015..2021)
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 43
$ws.Range("E10").Value = 38
$ws.Range("F10").Value = 45
$ws.Range("G10").Value = 80
$ws.Range("H10").Value = 121

# Update Total row (row 11) to reflect the new totals
$ws.Range("B11").Value = 216
$ws.Range("C11").Value = 415
$ws.Range("D11").Value = 594
$ws.Range("E11").Value = 528
$ws.Range("F11").Value = 394
$ws.Range("G11").Value = 864
$ws.Range("H11").Value = 1191
